$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.999610960483551
$ws.Range("B1").Value = 2.110814332962036
$ws.Range("C1").Value = 6.909176826477051
$ws.Range("D1").Value = 2.036569595336914
$ws.Range("E1").Value = 1.374456167221069
